# The deck currently has its active Design ("Integral") stored in the
# theme part used by the slide master (ppt/theme/theme2.xml), while the
# Notes Master still refers to the original, no-longer-applied
# "Office Theme" colours (ppt/theme/theme1.xml).
#
# This commit re-applies the stock "Office Theme" colour palette to the
# presentation's active design, which PowerPoint persists by rewriting
# the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# of the currently-applied theme part in place.
#
# PowerPoint's RGB colour values are standard Windows COLORREFs
# (0x00BBGGRR), i.e. byte-reversed relative to the RRGGBB hex strings
# used in the OOXML <a:srgbClr val="RRGGBB"/> markup.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x00000000   # dk1      -> 000000
$tcs.Colors(2).RGB  = 0x00FFFFFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x006A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0x00E6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0x00D59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x00317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0x00A5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x0000C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0x00C47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x0047AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0x00C16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x00724F95   # folHlink -> 954F72
